$wb = $excel.ActiveWorkbook

# Sheet "OFF" - Week's update to row 2 (H row)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 197
$wsOff.Range("C2").Value = 128
$wsOff.Range("D2").Value = 39
$wsOff.Range("E2").Value = 22
$wsOff.Range("G2").Value = 6

# Sheet "DEF" - Week's update to row 2 (H row)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 182
$wsDef.Range("C2").Value = 137
$wsDef.Range("D2").Value = 48
$wsDef.Range("E2").Value = 31
